# Fix typos in the "how-to-use" deck (slide 1):
#  - Table 4 (stiffness matrix): row2 had "K","0" swapped -> should read "0","K"
#  - Table 12 (damping matrix): the lone "B" in row3 should read "B2"
#  - Rectangle 24 (red box around mass 1) resized/repositioned
#  - Two new labels ("For spring" / "For damper") added near the bottom tables

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Fix stiffness-matrix table (shape named "Table 4") ---
$tbl1 = $s.Shapes.Item("Table 4").Table
$tbl1.Cell(2,1).Shape.TextFrame.TextRange.Text = "0"
$tbl1.Cell(2,2).Shape.TextFrame.TextRange.Text = "K"

# --- Fix damping-matrix table (shape named "Table 12") ---
$tbl2 = $s.Shapes.Item("Table 12").Table
$tbl2.Cell(3,2).Shape.TextFrame.TextRange.Text = "B2"

# PowerPoint COM works in points (1 pt = 12700 EMU)
$emu = 12700

# --- Reposition/resize the red box around mass 1 (shape "Rectangle 24") ---
$rect = $s.Shapes.Item("Rectangle 24")
$rect.Left = 251520 / $emu
$rect.Top = 2749570 / $emu
$rect.Width = 576064 / $emu
$rect.Height = 409982 / $emu

# --- Add "For spring" label ---
$tb1 = $s.Shapes.AddTextbox(1, 4174530 / $emu, 4509120 / $emu, 1117550 / $emu, 369332 / $emu)
$tb1.Name = "TextBox 1"
$tb1.TextFrame.WordWrap = 0
$tb1.TextFrame.TextRange.Text = "For spring"

# --- Add "For damper" label ---
$tb2 = $s.Shapes.AddTextbox(1, 6905528 / $emu, 4538230 / $emu, 1276247 / $emu, 369332 / $emu)
$tb2.Name = "TextBox 22"
$tb2.TextFrame.WordWrap = 0
$tb2.TextFrame.TextRange.Text = "For damper"
